$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species/record data between row 13 and row 14
# (columns A, B, E, F, G, H, P, Q, R, S), leaving all other columns as-is.

$ws.Range("A13").Value = 111703377
$ws.Range("B13").Value = 90678
$ws.Range("E13").Value = 4366
$ws.Range("F13").Value = "Skarp dropptaggsvamp"
$ws.Range("G13").Value = "Hydnellum peckii"
$ws.Range("H13").Value = "Banker"
$ws.Range("P13").Value = "N Björklunda (N Björklunda), Nrk"
$ws.Range("Q13").Value = 516440.1511331969
$ws.Range("R13").Value = 6574461.330051985
$ws.Range("S13").Value = 10

$ws.Range("A14").Value = 111704069
$ws.Range("B14").Value = 84997
$ws.Range("E14").Value = 3279
$ws.Range("F14").Value = "Maskfingersvamp"
$ws.Range("G14").Value = "Clavaria fragilis"
$ws.Range("H14").Value = "Holmsk.:Fr."
$ws.Range("P14").Value = "N Björklunda  (N Björklunda ), Nrk"
$ws.Range("Q14").Value = 516437.0004434386
$ws.Range("R14").Value = 6574480.208031038
$ws.Range("S14").Value = 15
